# Daily attendance processing - 2026-01-08 03:33:39
# Swap the order of "Recorded By" entries in column G from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}
